# Add a new "weatherForecast" localization row to the language sheet.
# Mirrors the row directly above it (row 30) for formatting (wrap-text style),
# then updates the active selection to the newly added cell, matching
# Excel's own behaviour when a user types a new row at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = "weatherForecast"
$ws.Range("B31").Value = "Weather Forecast"

# Copy formatting (wrap-text style) from the row above so the new row matches
# the rest of the table.
$ws.Range("B30").Copy()
$ws.Range("B31").PasteSpecial(-4122)

[void]$ws.Range("B31").Select()
